$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date value from 45170 to 45174 for rows 2-19
$ws.Range("C2:C19").Value = 45174
